$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# Add the new header "Correction" in N1, matching the style used by the
# other header cells (row 1 uses a bold/centered/bordered style).
$ws.Range("N1").Value = "Correction"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the previously "blank" M-column cells (rows 2-13) with the text "nan"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"
}

# Create blank *text* cells in the new N column (rows 2-13). A leading
# single-quote forces Excel to store the cell as (empty) text instead of
# clearing it; copying the plain formatting back in removes the quote-
# prefix style so the cell ends up styled like the rest of the data rows.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 14).Value = "'"
}
$ws.Range("A2").Copy()
$ws.Range("N2:N13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
